# Scan profiles plot changed from Bn to bn
# Update "Bn Roxie" columns (S:AF) in rows 2-8 with recalculated bn-based values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("S2").Value = [double]"0.0003374155165390692"
$ws.Range("T2").Value = [double]"0.001337716227943058"
$ws.Range("U2").Value = [double]"-3.836301416762852E-06"
$ws.Range("V2").Value = [double]"-7.935477442629211E-05"
$ws.Range("W2").Value = [double]"-9.125363393711146E-07"
$ws.Range("X2").Value = [double]"0.000103319795918088"
$ws.Range("Y2").Value = [double]"9.641793466235713E-07"
$ws.Range("Z2").Value = [double]"3.042224649662914E-05"
$ws.Range("AA2").Value = [double]"-2.557583858044626E-08"
$ws.Range("AB2").Value = [double]"4.90838083535676E-05"
$ws.Range("AC2").Value = [double]"-1.007132172838627E-06"
$ws.Range("AD2").Value = [double]"-0.0002135352256541946"
$ws.Range("AE2").Value = [double]"5.455179660721536E-07"
$ws.Range("AF2").Value = [double]"0.0001170906689813567"
# Row 3
$ws.Range("S3").Value = [double]"0.003145268273837382"
$ws.Range("T3").Value = [double]"-5.159703901125607"
$ws.Range("U3").Value = [double]"0.0002426066663963457"
$ws.Range("V3").Value = [double]"-3.447569149519026"
$ws.Range("W3").Value = [double]"-1.342709741750565E-05"
$ws.Range("X3").Value = [double]"-3.807632378352325"
$ws.Range("Y3").Value = [double]"1.803884600817608E-06"
$ws.Range("Z3").Value = [double]"-0.7959591041396851"
$ws.Range("AA3").Value = [double]"5.387616708110245E-07"
$ws.Range("AB3").Value = [double]"0.06023825116324097"
$ws.Range("AC3").Value = [double]"-1.362621105055171E-06"
$ws.Range("AD3").Value = [double]"-0.6113731744562918"
$ws.Range("AE3").Value = [double]"5.870551755424866E-07"
$ws.Range("AF3").Value = [double]"0.01608037360569016"
# Row 4
$ws.Range("S4").Value = [double]"-0.001687946010371525"
$ws.Range("T4").Value = [double]"-13.02621989082303"
$ws.Range("U4").Value = [double]"-0.0006945508100858341"
$ws.Range("V4").Value = [double]"0.2950884838163252"
$ws.Range("W4").Value = [double]"-4.516807969753136E-05"
$ws.Range("X4").Value = [double]"1.776289229535417"
$ws.Range("Y4").Value = [double]"3.556611376828276E-06"
$ws.Range("Z4").Value = [double]"2.579296971163119"
$ws.Range("AA4").Value = [double]"-5.966173124008846E-07"
$ws.Range("AB4").Value = [double]"3.086666735785285"
$ws.Range("AC4").Value = [double]"-1.791146375175698E-07"
$ws.Range("AD4").Value = [double]"-1.695279449585117"
$ws.Range("AE4").Value = [double]"-5.916188057279711E-07"
$ws.Range("AF4").Value = [double]"0.3237615017416037"
# Row 5
$ws.Range("AC5").Value = [double]"-5.406482465462275E-07"
# Row 6
$ws.Range("S6").Value = [double]"0.0008308706508825413"
$ws.Range("T6").Value = [double]"-12.99846623850372"
$ws.Range("U6").Value = [double]"-0.0002192530052636901"
$ws.Range("V6").Value = [double]"0.27769713835233"
$ws.Range("W6").Value = [double]"-4.419065289680839E-05"
$ws.Range("X6").Value = [double]"1.784511003164041"
$ws.Range("Y6").Value = [double]"-1.466536158456543E-06"
$ws.Range("Z6").Value = [double]"2.576825011735771"
$ws.Range("AA6").Value = [double]"-2.60310304880167E-06"
$ws.Range("AB6").Value = [double]"3.086920153558611"
$ws.Range("AC6").Value = [double]"9.151666291426928E-07"
$ws.Range("AD6").Value = [double]"-1.695043311474559"
$ws.Range("AE6").Value = [double]"-2.788866050663632E-06"
$ws.Range("AF6").Value = [double]"0.3235417131396416"
# Row 7
$ws.Range("S7").Value = [double]"-0.001703083995411835"
$ws.Range("T7").Value = [double]"-2.912983473620887"
$ws.Range("U7").Value = [double]"-0.000122347679660836"
$ws.Range("V7").Value = [double]"-0.07706840814960061"
$ws.Range("W7").Value = [double]"2.349075611653704E-05"
$ws.Range("X7").Value = [double]"-2.024201225339"
$ws.Range("Y7").Value = [double]"1.976982359299177E-06"
$ws.Range("Z7").Value = [double]"-0.4878575412004966"
$ws.Range("AA7").Value = [double]"-1.790028598755125E-07"
$ws.Range("AB7").Value = [double]"0.2602097271581776"
$ws.Range("AC7").Value = [double]"4.500950413440616E-07"
$ws.Range("AD7").Value = [double]"-0.5309228916190798"
$ws.Range("AE7").Value = [double]"2.304550423522201E-06"
$ws.Range("AF7").Value = [double]"0.04986165914984727"
# Row 8
$ws.Range("S8").Value = [double]"-3.101201890367674E-05"
$ws.Range("T8").Value = [double]"0.004891586602117523"
$ws.Range("U8").Value = [double]"-2.079313957908111E-06"
$ws.Range("V8").Value = [double]"-2.278909536938637E-05"
$ws.Range("W8").Value = [double]"1.477197993531028E-06"
$ws.Range("X8").Value = [double]"9.616244135084262E-05"
$ws.Range("Y8").Value = [double]"7.749681568662174E-07"
$ws.Range("Z8").Value = [double]"2.836587273695381E-05"
$ws.Range("AA8").Value = [double]"2.19672190335622E-07"
$ws.Range("AB8").Value = [double]"5.50600451644013E-05"
$ws.Range("AC8").Value = [double]"-1.315272239759572E-06"
$ws.Range("AD8").Value = [double]"-0.0001700515840499074"
$ws.Range("AE8").Value = [double]"-8.665636477264422E-07"
$ws.Range("AF8").Value = [double]"6.589463835234729E-05"
